$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 108: new episode C2E107 (was a blank/gap row)
$ws.Cells.Item(108, 1).Value = "C2E107"
$ws.Cells.Item(108, 2).Value = 0
$ws.Cells.Item(108, 3).Value = 1134
$ws.Cells.Item(108, 4).Value = 0
$ws.Cells.Item(108, 5).Value = 0
$ws.Cells.Item(108, 6).Value = 1134
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = -4334
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = -4334
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = -3200
$ws.Cells.Item(108, 14).Value = 0
$ws.Cells.Item(108, 15).Value = 0
$ws.Cells.Item(108, 16).Value = -3200

# Row 109: previously the TOTALS summary row, now replaced with episode C2E108 data
$ws.Cells.Item(109, 1).Value = "C2E108"
$ws.Cells.Item(109, 2).Value = 0
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 0
$ws.Cells.Item(109, 5).Value = 0
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = -50
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = -50
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -50
$ws.Cells.Item(109, 14).Value = 0
$ws.Cells.Item(109, 15).Value = 0
$ws.Cells.Item(109, 16).Value = -50

# Row 110: episode C2E109
$ws.Cells.Item(110, 1).Value = "C2E109"
$ws.Cells.Item(110, 2).Value = 0
$ws.Cells.Item(110, 3).Value = 2400
$ws.Cells.Item(110, 4).Value = 0
$ws.Cells.Item(110, 5).Value = 0
$ws.Cells.Item(110, 6).Value = 2400
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = -16162
$ws.Cells.Item(110, 9).Value = -19
$ws.Cells.Item(110, 10).Value = -6
$ws.Cells.Item(110, 11).Value = -16163.96
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = -13762
$ws.Cells.Item(110, 14).Value = -19
$ws.Cells.Item(110, 15).Value = -6
$ws.Cells.Item(110, 16).Value = -13763.96

# Row 111: episode C2E110
$ws.Cells.Item(111, 1).Value = "C2E110"
$ws.Cells.Item(111, 2).Value = 0
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 0
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = -420
$ws.Cells.Item(111, 9).Value = -3
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = -420.3
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = -420
$ws.Cells.Item(111, 14).Value = -3
$ws.Cells.Item(111, 15).Value = 0
$ws.Cells.Item(111, 16).Value = -420.3

# Row 112: episode C2E111
$ws.Cells.Item(112, 1).Value = "C2E111"
$ws.Cells.Item(112, 2).Value = 546
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 0
$ws.Cells.Item(112, 6).Value = 5460
$ws.Cells.Item(112, 7).Value = -548
$ws.Cells.Item(112, 8).Value = -9961
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = -15441
$ws.Cells.Item(112, 12).Value = -2
$ws.Cells.Item(112, 13).Value = -9961
$ws.Cells.Item(112, 14).Value = 0
$ws.Cells.Item(112, 15).Value = 0
$ws.Cells.Item(112, 16).Value = -9981

# Row 113: episode C2E112
$ws.Cells.Item(113, 1).Value = "C2E112"
$ws.Cells.Item(113, 2).Value = 0
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 0
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = -607
$ws.Cells.Item(113, 9).Value = -37
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = -610.70000000000005
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -607
$ws.Cells.Item(113, 14).Value = -37
$ws.Cells.Item(113, 15).Value = 0
$ws.Cells.Item(113, 16).Value = -610.70000000000005

# Row 114: episode C2E113
$ws.Cells.Item(114, 1).Value = "C2E113"
$ws.Cells.Item(114, 2).Value = 0
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 0
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 0
$ws.Cells.Item(114, 14).Value = 0
$ws.Cells.Item(114, 15).Value = 0
$ws.Cells.Item(114, 16).Value = 0

# Update the selected cell to match the final workbook state
$ws.Range("Q104").Select() | Out-Null

Write-Output "done"
